# Feria Lagunitas de Puerto Montt - Brócoli
# Insert two new weekly records at the top of the data block (rows 234-235),
# pushing the existing rows 234..329 down to 236..331.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 234 (shifts old 234..329 -> 236..331)
$ws.Range("A234:A235").EntireRow.Insert()

# New row 234: Primera calidad, fecha 44704 (2022-05-23), Región del Maule
$ws.Range("A234").Value = 4
$ws.Range("B234").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C234").Value = "Los Lagos"
$ws.Range("D234").Value = 44704
$ws.Range("E234").Value = 10
$ws.Range("F234").Value = 100112023
$ws.Range("G234").Value = "Brócoli"
$ws.Range("H234").Value = "Sin especificar"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 250
$ws.Range("K234").Value = 1500
$ws.Range("L234").Value = 1500
$ws.Range("M234").Value = 1500
$ws.Range("N234").Value = "$/unidad"
$ws.Range("O234").Value = "Región del Maule"
$ws.Range("P234").Value = 1500
$ws.Range("Q234").Value = 1
$ws.Range("R234").Value = "Hortaliza"

# New row 235: Segunda calidad, fecha 44704 (2022-05-23), Región del Maule
$ws.Range("A235").Value = 4
$ws.Range("B235").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C235").Value = "Los Lagos"
$ws.Range("D235").Value = 44704
$ws.Range("E235").Value = 10
$ws.Range("F235").Value = 100112023
$ws.Range("G235").Value = "Brócoli"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Segunda"
$ws.Range("J235").Value = 250
$ws.Range("K235").Value = 1200
$ws.Range("L235").Value = 1200
$ws.Range("M235").Value = 1200
$ws.Range("N235").Value = "$/unidad"
$ws.Range("O235").Value = "Región del Maule"
$ws.Range("P235").Value = 1200
$ws.Range("Q235").Value = 1
$ws.Range("R235").Value = "Hortaliza"
